# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The detail table (rows 16-27, columns C:G) is rebuilt: instead of
# alternating between the two workers period by period, the rows are
# regrouped so all of JORGE LUIS MIRANDA RAMOS's periods come first
# (newest period 1906 down to oldest 1901), followed by all of
# LINA MARIA MAZO MONSALVE's periods (also 1906 down to 1901). The
# "Valor Mora" (F) amounts travel with their owning worker/period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nit1  = "73270932"
$name1 = "JORGE LUIS MIRANDA RAMOS"
$nit2  = "43600251"
$name2 = "LINA MARIA MAZO MONSALVE"

# row => (NIT, Nombre, Periodo, Valor Mora, Salario Basico)
$rows = @{
  16 = @($nit1, $name1, "1906", 22666,  1000000)
  17 = @($nit1, $name1, "1905", 40000,  1000000)
  18 = @($nit1, $name1, "1904", 40000,  1000000)
  19 = @($nit1, $name1, "1903", 40000,  1000000)
  20 = @($nit1, $name1, "1902", 40000,  1000000)
  21 = @($nit1, $name1, "1901", 40000,  1000000)
  22 = @($nit2, $name2, "1906", 56666,  2500000)
  23 = @($nit2, $name2, "1905", 100000, 2500000)
  24 = @($nit2, $name2, "1904", 100000, 2500000)
  25 = @($nit2, $name2, "1903", 100000, 2500000)
  26 = @($nit2, $name2, "1902", 100000, 2500000)
  27 = @($nit2, $name2, "1901", 100000, 2500000)
}

foreach ($r in 16..27) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]   # C - N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $vals[1]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $vals[2]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $vals[3]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G - Salario Basico
}
